{"js": "// The MoU template had a stray \"referral\" placeholder block spread across\n// three paragraphs:\n//   1) \"${referral}\" + \"{{\"\n//   2) \"${referral_signature:100px:100px:ratio=true}\"\n//   3) \"${referral_name}\" + \"email}}\"\n// The fix collapses this down to a single \"}}\" run (removing the whole\n// referral placeholder) in the first paragraph, and removes the other two\n// paragraphs entirely.\n\nconst signatureMatches = context.document.body.search(\"${referral_signature\", { matchCase: true });\nconst nameMatches = context.document.body.search(\"${referral_name}\", { matchCase: true });\nsignatureMatches.load(\"items\");\nnameMatches.load(\"items\");\nawait context.sync();\n\nif (signatureMatches.items.length > 0 && nameMatches.items.length > 0) {\n  // Remove the two now-obsolete paragraphs entirely (including their\n  // paragraph marks/properties).\n  nameMatches.items[0].paragraphs.getFirst().delete();\n  signatureMatches.items[0].paragraphs.getFirst().delete();\n  await context.sync();\n\n  // Drop the leading \"${referral}\" run in the remaining paragraph.\n  const referralMatches = context.document.body.search(\"${referral}\", { matchCase: true });\n  referralMatches.load(\"items\");\n  await context.sync();\n  if (referralMatches.items.length > 0) {\n    referralMatches.items[0].delete();\n    await context.sync();\n  }\n\n  // Turn the surviving \"{{\" run into \"}}\".\n  const braceMatches = context.document.body.search(\"{{\", { matchCase: true });\n  braceMatches.load(\"items\");\n  await context.sync();\n  if (braceMatches.items.length > 0) {\n    braceMatches.items[0].insertText(\"}}\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# The MoU template had a stray \"referral\" placeholder block spread across\n# three paragraphs:\n#   1) \"${referral}\" + \"{{\"\n#   2) \"${referral_signature:100px:100px:ratio=true}\"\n#   3) \"${referral_name}\" + \"email}}\"\n# The fix collapses this down to a single \"}}\" run (removing the whole\n# referral placeholder) in the first paragraph, and removes the other two\n# paragraphs entirely.\n\n$idxFirst = -1\n$idxSignature = -1\n$idxName = -1\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.StartsWith('${referral}{{')) {\n        $idxFirst = $i\n    } elseif ($t.StartsWith('${referral_signature')) {\n        $idxSignature = $i\n    } elseif ($t.StartsWith('${referral_name}')) {\n        $idxName = $i\n    }\n}\n\nif ($idxFirst -eq -1 -or $idxSignature -eq -1 -or $idxName -eq -1) {\n    Write-Output \"Target referral paragraphs not found; no changes made.\"\n} else {\n    # Remove the trailing two paragraphs completely (highest index first so\n    # the lower indices stay valid while we work).\n    $d.Paragraphs.Item($idxName).Range.Delete()\n    $d.Paragraphs.Item($idxSignature).Range.Delete()\n\n    # In the remaining paragraph, drop the leading \"${referral}\" run and\n    # turn the surviving \"{{\" run into \"}}\".\n    $p = $d.Paragraphs.Item($idxFirst)\n    $prefix = '${referral}'\n    $leadRange = $d.Range($p.Range.Start, $p.Range.Start + $prefix.Length)\n    $leadRange.Delete()\n    $p.Range.Text = \"}}\"\n\n    Write-Output \"Removed referral placeholder block.\"\n}\n"}
